$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheetId=1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1592
$ws1.Range("F4").Value = 52

# Update "全部类型" sheet (sheetId=4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1592
$ws4.Range("F4").Value = 52
